$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2799.8
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 3124.75
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 3124.75
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -3776.75
$ws.Range("H39").Value = 584.1739
$ws.Range("I39").Value = 607.53845
$ws.Range("J39").Value = 553.8
$ws.Range("K39").Value = 1822.61535
$ws.Range("L39").Value = 1661.4
$ws.Range("M39").Value = -1526.61535
$ws.Range("N39").Value = -2253.4
$ws.Range("H43").Value = 4630243.5
$ws.Range("I43").Value = 495.375
$ws.Range("J43").Value = 13889740
$ws.Range("K43").Value = 495.375
$ws.Range("L43").Value = 13889740
$ws.Range("M43").Value = -426.375
$ws.Range("N43").Value = -13889878
$ws.Range("H69").Value = 3973.3333
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 3973.3333
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H97").Value = 2217.6843
$ws.Range("J97").Value = 2217.6843
$ws.Range("L97").Value = 6653.0529
$ws.Range("N97").Value = -7645.0529
$ws.Range("H125").Value = 5000
$ws.Range("J125").Value = 5000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -49920
$ws.Range("H132").Value = 8775714
$ws.Range("I132").Value = 12822569
$ws.Range("K132").Value = 38467707
$ws.Range("M132").Value = -38465177
$ws.Range("H138").Value = 2178.0103
$ws.Range("J138").Value = 2283.3604
$ws.Range("L138").Value = 6850.081200000001
$ws.Range("N138").Value = -17130.0812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2958.875
$ws.Range("I32").Value = 3110.3777
$ws.Range("K32").Value = 3110.3777
$ws.Range("M32").Value = -2823.3777
$ws.Range("H41").Value = 4859.4287
$ws.Range("I41").Value = 4609
$ws.Range("J41").Value = 6362
$ws.Range("K41").Value = 4609
$ws.Range("L41").Value = 6362
$ws.Range("M41").Value = -4195
$ws.Range("N41").Value = -7190
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H102").Value = 27795046
$ws.Range("I102").Value = 27795046
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 27795046
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -27793424
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 71431390
$ws.Range("I105").Value = 125002690
$ws.Range("K105").Value = 125002690
$ws.Range("M105").Value = -125000943
$ws.Range("H107").Value = 1633.3846
$ws.Range("I107").Value = 1302.625
$ws.Range("K107").Value = 1302.625
$ws.Range("M107").Value = 617.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1814.4445
$ws.Range("I5").Value = 971.6667
$ws.Range("J5").Value = 3500
$ws.Range("K5").Value = 971.6667
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = -859.6667
$ws.Range("N5").Value = -3724
$ws.Range("H31").Value = 940.2245
$ws.Range("I31").Value = 801.2286
$ws.Range("K31").Value = 801.2286
$ws.Range("M31").Value = -506.2286
$ws.Range("H34").Value = 940.2245
$ws.Range("I34").Value = 801.2286
$ws.Range("K34").Value = 801.2286
$ws.Range("M34").Value = -599.2286
$ws.Range("H99").Value = 2075.077
$ws.Range("I99").Value = 1883.1111
$ws.Range("K99").Value = 1883.1111
$ws.Range("M99").Value = -385.1111000000001
$ws.Range("H122").Value = 764.8
$ws.Range("J122").Value = 1200
$ws.Range("L122").Value = 3600
$ws.Range("N122").Value = -8500
$ws.Range("H126").Value = 2075.077
$ws.Range("I126").Value = 1883.1111
$ws.Range("K126").Value = 5649.3333
$ws.Range("M126").Value = -3179.3333
$ws.Range("H134").Value = 8131236.5
$ws.Range("I134").Value = 9010047
$ws.Range("J134").Value = 2239.5
$ws.Range("K134").Value = 27030141
$ws.Range("L134").Value = 6718.5
$ws.Range("M134").Value = -27027606
$ws.Range("N134").Value = -11788.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1382.8485
$ws.Range("I68").Value = 699.61536
$ws.Range("K68").Value = 2098.84608
$ws.Range("M68").Value = -1287.84608
$ws.Range("H71").Value = 1382.8485
$ws.Range("I71").Value = 699.61536
$ws.Range("K71").Value = 6296.53824
$ws.Range("M71").Value = -2240.53824
$ws.Range("H88").Value = 1580
$ws.Range("J88").Value = 6800
$ws.Range("L88").Value = 20400
$ws.Range("N88").Value = -21256
$ws.Range("H91").Value = 1580
$ws.Range("J91").Value = 6800
$ws.Range("L91").Value = 20400
$ws.Range("N91").Value = -23364
$ws.Range("H103").Value = 2966.75
$ws.Range("I103").Value = 1134.6
$ws.Range("J103").Value = 4275.4287
$ws.Range("K103").Value = 3403.8
$ws.Range("L103").Value = 12826.2861
$ws.Range("M103").Value = -2524.8
$ws.Range("N103").Value = -14584.2861
$ws.Range("H126").Value = 5680
$ws.Range("J126").Value = 5680
$ws.Range("L126").Value = 17040
$ws.Range("N126").Value = -26920
$ws.Range("H131").Value = 17242654
$ws.Range("I131").Value = 142857870
$ws.Range("J131").Value = 1350.2745
$ws.Range("K131").Value = 428573610
$ws.Range("L131").Value = 4050.8235
$ws.Range("M131").Value = -428568570
$ws.Range("N131").Value = -14130.8235
$ws.Range("H137").Value = 9953
$ws.Range("J137").Value = 11414.777
$ws.Range("L137").Value = 34244.331
$ws.Range("N137").Value = -44444.331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 20250
$ws.Range("J46").Value = 20250
$ws.Range("L46").Value = 20250
$ws.Range("N46").Value = -20562
$ws.Range("H52").Value = 19980
$ws.Range("J52").Value = 19980
$ws.Range("L52").Value = 19980
$ws.Range("N52").Value = -20498
$ws.Range("H57").Value = 19999.9
$ws.Range("J57").Value = 19999.9
$ws.Range("L57").Value = 19999.9
$ws.Range("N57").Value = -21639.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2051.4443
$ws.Range("I7").Value = 1932.25
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 1932.25
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -1820.25
$ws.Range("N7").Value = -3229
$ws.Range("H16").Value = 1158.4375
$ws.Range("I16").Value = 1048.6923
$ws.Range("J16").Value = 1634
$ws.Range("K16").Value = 1048.6923
$ws.Range("L16").Value = 1634
$ws.Range("M16").Value = -878.6922999999999
$ws.Range("N16").Value = -1974
$ws.Range("H38").Value = 4900
$ws.Range("I38").Value = 4900
$ws.Range("K38").Value = 4900
$ws.Range("M38").Value = -4490
$ws.Range("H40").Value = 2261.2307
$ws.Range("I40").Value = 1988.7
$ws.Range("K40").Value = 1988.7
$ws.Range("M40").Value = -1852.7
$ws.Range("H126").Value = 2051.4443
$ws.Range("I126").Value = 1932.25
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 5796.75
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -3326.75
$ws.Range("N126").Value = -13955
$ws.Range("H136").Value = 1655.5834
$ws.Range("I136").Value = 1556.2
$ws.Range("K136").Value = 4668.6
$ws.Range("M136").Value = -2118.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 24500
$ws.Range("J94").Value = 24500
$ws.Range("L94").Value = 24500
$ws.Range("N94").Value = -26302
$ws.Range("H126").Value = 66667972
$ws.Range("I126").Value = 111111740
$ws.Range("J126").Value = 2316.5
$ws.Range("K126").Value = 333335220
$ws.Range("L126").Value = 6949.5
$ws.Range("M126").Value = -333332750
$ws.Range("N126").Value = -11889.5
$ws.Range("H136").Value = 1806
$ws.Range("I136").Value = 1649.1
$ws.Range("K136").Value = 4947.299999999999
$ws.Range("M136").Value = -2397.299999999999
